$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 40; $r -le 52; $r++) {
    $ws.Range("O$r").Value = "['Costa Rica', 'Ireland']"
}

for ($r = 53; $r -le 59; $r++) {
    $ws.Range("O$r").Value = "['Costa Rica', 'Argentina', 'Colombia', 'Ireland']"
}

for ($r = 63; $r -le 73; $r++) {
    $ws.Range("O$r").Value = "['Argentina', 'Colombia', 'Scotland', 'Austria']"
}
